# Generate Report for Handoff
# Updates the localization status report: refreshes the generated/handoff
# timestamps and bumps priority from "low" to "ht" for the rows that
# correspond to the 0c7147df... / 1b5f9af8... / 8ac4a81b... / 9165827b...
# files (rows 4-7) on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# New timestamps produced by the report generation run.
$newHoGenerateDate = "2016-08-12 14:34:25"
$newZhCnHandoffDate = "2016-08-12 14:34:18"

# Overview sheet: "Latest HO Xliff Generate Date" column G, rows 4-7.
$overview.Range("G4:G7").Value = $newHoGenerateDate

# zh-cn sheet: Priority column E rows 4-7 moves from "low" to "ht";
# "Latest Handoff Datetime" column H rows 4-7 refreshed.
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = $newZhCnHandoffDate

# de-de sheet: Priority column E rows 4-7 moves from "low" to "ht";
# "Latest Handoff Datetime" column H rows 4-7 shares the HO generate date.
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = $newHoGenerateDate
